$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
# Force a real text-replacement so the runtime actually rebuilds the runs
# (setting directly to the already-matching concatenated text is a no-op).
$sh.TextFrame.TextRange.Text = "placeholder"
$sh.TextFrame.TextRange.Text = "The picture first"
